# Insert a new "column_start" parameter row into the "General parameters"
# sheet (between "year_end" and "rho"), update the defined names that point
# at rows shifted down by the insert, and refresh the view selections to
# reflect the state after the edit (matches the author's commit: "allow the
# user to run tests on a flexible part of the entered time series to speed
# up simulation time").

$wb = $excel.ActiveWorkbook

# --- "General parameters" sheet: insert new row 3 -------------------------
$wsGeneral = $wb.Worksheets.Item("General parameters")

# Shift existing rows 3:7 (rho, g, T_fill_thres, LOEE_allowed, option_storage)
# down to rows 4:8, freeing up row 3 for the new parameter.
$wsGeneral.Rows("3:3").Insert()

# Populate the newly freed row 3 with the new parameter.
$wsGeneral.Range("A3").Value = "column_start"
$wsGeneral.Range("B3").Value = 1
$wsGeneral.Range("C3").Value = "index of column (first column = 1) corresponding to year_start in time series Excel sheets (this needs to be the same across all Excel sheets)"

# Match the formatting used by the other "note" cells in column C (wrapped
# text), and size the row to fit the (longer) wrapped description.
$wsGeneral.Range("C3").WrapText = $true
$wsGeneral.Rows("3:3").RowHeight = 43.5

# --- Workbook-level defined names: re-point at the shifted rows -----------
# (Delete + re-Add rather than assigning .RefersTo so the dependent
# formulas on "Hydropower plant parameters" - which use rho/g - get
# re-resolved against the new targets.)
$wb.Names.Item("rho").Delete()
$wb.Names.Add("rho", "='General parameters'!`$B`$4")
$wb.Names.Item("g").Delete()
$wb.Names.Add("g", "='General parameters'!`$B`$5")
$wb.Names.Item("T_fill_thres").Delete()
$wb.Names.Add("T_fill_thres", "='General parameters'!`$B`$6")
$wb.Names.Item("LOEE_allowed").Delete()
$wb.Names.Add("LOEE_allowed", "='General parameters'!`$B`$7")
$wb.Names.Item("option_storage").Delete()
$wb.Names.Add("option_storage", "='General parameters'!`$B`$8")

# --- Update selections on both touched sheets ------------------------------
$wsHydro = $wb.Worksheets.Item("Hydropower plant parameters")
[void]$wsHydro.Range("B16").Select()

[void]$wsGeneral.Range("C2").Select()
$wsGeneral.Activate()

Write-Output "done"
